$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Country"
$ws.Range("C2").Value = "Peru"
$ws.Range("C3").Value = "Peru"
$ws.Range("C4").Value = "Peru"

$ws.Range("D1").Value = "Region"
$ws.Range("D2").Value = "east"
$ws.Range("D3").Value = "east"
$ws.Range("D4").Value = "west"

$ws.Range("E1").Value = "Line of Business"
$ws.Range("E2").Value = "Construction"
$ws.Range("E3").Value = "Construction"
$ws.Range("E4").Value = "Construction"

$ws.Columns.Item(2).ColumnWidth = 24.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.5

$ws.Range("G5").Select()
